# Generate Report for Handoff
# Adds a new handed-off file ("e72579d7-37cd-46dc-a079-1ddb3a90ee3b.md") as
# row 9 of the Overview / zh-cn / de-de tables.

$wb = $excel.ActiveWorkbook

$baseName   = "e72579d7-37cd-46dc-a079-1ddb3a90ee3b.md"
$pathName   = "e2e\" + $baseName
$xliffZh    = "e72579d7-37cd-46dc-a079-1ddb3a90ee3b.a731d6e910daae84f90505762c81f6c56c639889.zh-cn.xlf"
$xliffDe    = "e72579d7-37cd-46dc-a079-1ddb3a90ee3b.a731d6e910daae84f90505762c81f6c56c639889.de-de.xlf"
$dateZh     = "2016-08-17 00:41:20"
$dateDe     = "2016-08-17 00:41:25"
$hrefUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a731d6e910daae84f90505762c81f6c56c639889/e2e/$baseName"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 9
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $baseName
$wsOverview.Range("B9").Value = $pathName
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = $dateDe
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $hrefUrl, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 9
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A9").Value = $baseName
$wsZh.Range("B9").Value = ".md"
$wsZh.Range("C9").Value = "Ready for handoff"
$wsZh.Range("D9").Value = "e2e"
$wsZh.Range("E9").Value = "ht"
$wsZh.Range("F9").Value = "False"
$wsZh.Range("G9").Value = $xliffZh
$wsZh.Range("H9").Value = $dateZh
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I9").Value = ""
$wsZh.Range("J9").Value = ""
$wsZh.Range("K9").Value = "0001-01-01 00:00:00"
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L9").Value = ""
$wsZh.Range("M9").Value = "True"
$wsZh.Range("N9").Value = ""
$wsZh.Range("O9").Value = "False"
$wsZh.Range("P9").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), $hrefUrl, "", "", $baseName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 9
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A9").Value = $baseName
$wsDe.Range("B9").Value = ".md"
$wsDe.Range("C9").Value = "Ready for handoff"
$wsDe.Range("D9").Value = "e2e"
$wsDe.Range("E9").Value = "ht"
$wsDe.Range("F9").Value = "False"
$wsDe.Range("G9").Value = $xliffDe
$wsDe.Range("H9").Value = $dateDe
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I9").Value = ""
$wsDe.Range("J9").Value = ""
$wsDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L9").Value = ""
$wsDe.Range("M9").Value = "True"
$wsDe.Range("N9").Value = ""
$wsDe.Range("O9").Value = "False"
$wsDe.Range("P9").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), $hrefUrl, "", "", $baseName) | Out-Null
